# Replace the single paragraph in the header/footer with a 2x3 table,
# matching the structure:
#   row1: a | b | c
#   row2: d | <picture> | "Page X of Y" field codes
# (identical table layout is used for both the header and the footer).

$tableXml = @'
<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:w10="urn:schemas-microsoft-com:office:word">
    <w:tblGrid>
      <w:gridCol w:w="1" w:type="dxa"/>
      <w:gridCol w:w="1" w:type="dxa"/>
      <w:gridCol w:w="1" w:type="dxa"/>
    </w:tblGrid>
    <w:tr>
      <w:trPr/>
      <w:tc>
        <w:tcPr>
          <w:tcW w:w="1" w:type="dxa"/>
        </w:tcPr>
        <w:p>
          <w:pPr/>
          <w:r>
            <w:rPr/>
            <w:t xml:space="preserve">a</w:t>
          </w:r>
        </w:p>
      </w:tc>
      <w:tc>
        <w:tcPr>
          <w:tcW w:w="1" w:type="dxa"/>
        </w:tcPr>
        <w:p>
          <w:pPr/>
          <w:r>
            <w:rPr/>
            <w:t xml:space="preserve">b</w:t>
          </w:r>
        </w:p>
      </w:tc>
      <w:tc>
        <w:tcPr>
          <w:tcW w:w="1" w:type="dxa"/>
        </w:tcPr>
        <w:p>
          <w:pPr/>
          <w:r>
            <w:rPr/>
            <w:t xml:space="preserve">c</w:t>
          </w:r>
        </w:p>
      </w:tc>
    </w:tr>
    <w:tr>
      <w:trPr/>
      <w:tc>
        <w:tcPr>
          <w:tcW w:w="1" w:type="dxa"/>
        </w:tcPr>
        <w:p>
          <w:pPr/>
          <w:r>
            <w:rPr/>
            <w:t xml:space="preserve">d</w:t>
          </w:r>
        </w:p>
      </w:tc>
      <w:tc>
        <w:tcPr>
          <w:tcW w:w="1" w:type="dxa"/>
        </w:tcPr>
        <w:p>
          <w:r>
            <w:pict>
              <v:shape type="#_x0000_t75" style="width:80pt; height:80pt; margin-left:0pt; margin-top:0pt; mso-position-horizontal:left; mso-position-vertical:top; mso-position-horizontal-relative:char; mso-position-vertical-relative:line;">
                <w10:wrap type="inline"/>
                <v:imagedata r:id="rId1" o:title=""/>
              </v:shape>
            </w:pict>
          </w:r>
        </w:p>
      </w:tc>
      <w:tc>
        <w:tcPr>
          <w:tcW w:w="1" w:type="dxa"/>
        </w:tcPr>
        <w:p>
          <w:r>
            <w:t xml:space="preserve">Page </w:t>
          </w:r>
          <w:r>
            <w:fldChar w:fldCharType="begin"/>
          </w:r>
          <w:r>
            <w:instrText xml:space="preserve">PAGE</w:instrText>
          </w:r>
          <w:r>
            <w:fldChar w:fldCharType="separate"/>
          </w:r>
          <w:r>
            <w:fldChar w:fldCharType="end"/>
          </w:r>
          <w:r>
            <w:t xml:space="preserve"> of </w:t>
          </w:r>
          <w:r>
            <w:fldChar w:fldCharType="begin"/>
          </w:r>
          <w:r>
            <w:instrText xml:space="preserve">NUMPAGES</w:instrText>
          </w:r>
          <w:r>
            <w:fldChar w:fldCharType="separate"/>
          </w:r>
          <w:r>
            <w:fldChar w:fldCharType="end"/>
          </w:r>
          <w:r>
            <w:t xml:space="preserve">.</w:t>
          </w:r>
        </w:p>
      </w:tc>
    </w:tr>
  </w:tbl>
'@

$d = $word.ActiveDocument
$section = $d.Sections.Item(1)

$header = $section.Headers.Item(1)
$header.Range.InsertXML($tableXml)

$footer = $section.Footers.Item(1)
$footer.Range.InsertXML($tableXml)

Write-Output "header/footer tables inserted"
